$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price / volume figures scraped on Thu Apr 20 16:47:33 UTC 2023

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.917.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.971.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -4.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4084"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08648"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.068"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.961.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.755"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.288"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.013"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001074"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06635"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.011"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.816"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.932.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.183.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.063"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.182"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("E32").Value = "  -4.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09691"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.479"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.736"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.701"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02367"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.285"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06284"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.876"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.012"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -7.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.348"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6015"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.455"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000339"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.150"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.77%  "
